$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, A, B, C, D -- ordered to match the original authoring sequence
# so the shared-string table is rebuilt in the same order as the source workbook.
$data = @(
    ,@(381, "Button", 25052, "Kunder", "Lås upp Inaktivera")
    ,@(382, "Button", 23978, "Kunder", "Lås upp Landskod")
    ,@(383, "Button", 23983, "Kunder", "Lås upp Landskod leveransadress")
    ,@(384, "Button", 23868, "Kunder", "Lås upp Kundkategori")
    ,@(385, "Button", 23869, "Kunder", "Lås upp Distrikt")
    ,@(386, "Button", 23870, "Kunder", "Lås upp Säljare")
    ,@(387, "Button", 23873, "Kunder", "Lås up Leveranssätt")
    ,@(388, "Button", 23874, "Kunder", "Lås upp Speditör")
    ,@(389, "Button", 23875, "Kunder", "Lås upp Språk")
    ,@(390, "Button", 23876, "Kunder", "Lås upp Valuta")
    ,@(391, "Button", 23877, "Kunder", "Lås upp Sorteringsbegrepp")
    ,@(392, "Button", 23878, "Kunder", "Lås upp Kreditgräns kronor")
    ,@(393, "Button", 23879, "Kunder", "Lås upp Export")
    ,@(394, "Button", 23880, "Kunder", "Lås upp EU-kund")
    ,@(395, "Button", 23881, "Kunder", "Lås upp Räntefakturering")
    ,@(396, "Button", 23882, "Kunder", "Lås upp Betalningspåminnelse")
    ,@(419, "Button", 24879, "Kunder", "Lås upp Betalningspåminnelse")
    ,@(397, "Button", 23883, "Kunder", "Lås upp Påminnelseavgift")
    ,@(398, "Button", 23884, "Kunder", "Lås upp Restnotera ej")
    ,@(399, "Button", 23885, "Kunder", "Lås upp Expeditionsavgift")
    ,@(400, "Button", 23886, "Kunder", "Lås upp Frakt")
    ,@(401, "Button", 26659, "Kunder", "Lås upp Intrastat")
    ,@(402, "Button", 23887, "Kunder", "Lås upp Samlingsfakturering")
    ,@(403, "Button", 24123, "Kunder", "Lås upp Överför adress till beställning")
    ,@(404, "Button", 26202, "Kunder", "Lås upp Anmärkning 1")
    ,@(405, "Button", 26203, "Kunder", "Lås upp Anmärkning 2")
    ,@(406, "Button", 23888, "Kunder", "Lås upp Rabattavtal")
    ,@(407, "Button", 23871, "Kunder", "Lås upp Betalningsvillkor")
    ,@(410, "Button", 23889, "Kunder", "Lås upp Prislista")
    ,@(411, "Button", 23890, "Kunder", "Lås upp Fakturarabatt i %")
    ,@(412, "Button", 23891, "Kunder", "Lås upp Radrabatt")
    ,@(408, "Button", 23872, "Kunder", "Lås upp Leveransvillkor")
    ,@(409, "Button", 24088, "Kunder", "Lås upp Kundfordranskonto")
    ,@(413, "Button", 24873, "Kunder", "Lås upp Offert")
    ,@(414, "Button", 24874, "Kunder", "Lås upp Order")
    ,@(415, "Button", 24875, "Kunder", "Lås upp Följesedel")
    ,@(416, "Button", 24876, "Kunder", "Lås upp Extra orderdokument")
    ,@(417, "Button", 24877, "Kunder", "Lås upp Faktura")
    ,@(418, "Button", 24878, "Kunder", "Lås upp Extra fakturadokument")
    ,@(429, "Button", 24090, "Kunder", "Lås upp Kopia Offert")
    ,@(430, "Edit", 25007, "Kunder", "Kopia Offert")
    ,@(431, "Button", 23988, "Kunder", "Lås upp Kopia Order")
    ,@(432, "Edit", 24038, "Kunder", "Kopia Order")
    ,@(433, "Button", 23987, "Kunder", "Lås upp Kopia Följesedel")
    ,@(434, "Edit", 25003, "Kunder", "Kopia Följesedel")
    ,@(435, "Button", 24085, "Kunder", "Lås upp Kopia Extra orderdokument")
    ,@(436, "Edit", 25004, "Kunder", "Kopia Extra orderdokument")
    ,@(437, "Button", 23989, "Kunder", "Lås upp Kopia Faktura")
    ,@(438, "Edit", 24032, "Kunder", "Kopia Faktura")
    ,@(439, "Button", 23990, "Kunder", "Lås upp Kopia Kreditfaktura")
    ,@(440, "Edit", 24034, "Kunder", "Kopia Kreditfaktura")
    ,@(441, "Button", 23991, "Kunder", "Lås upp Kopia Räntefaktura")
    ,@(442, "Edit", 24036, "Kunder", "Kopia Räntefaktura")
    ,@(443, "Button", 24086, "Kunder", "Lås upp Kopia Kontantnota")
    ,@(444, "Edit", 25005, "Kunder", "Kopia Kontantnota")
    ,@(445, "Button", 24087, "Kunder", "Lås upp Kopia Extra fakturadokument")
    ,@(446, "Edit", 25006, "Kunder", "Kopia Extra fakturadokument")
    ,@(447, "Button", 23986, "Kunder", "Lås upp Dokumentmall Faktura")
    ,@(448, "Button", 24092, "Kunder", "Lås upp Dokumentmall Kreditfaktura")
    ,@(449, "Button", 24093, "Kunder", "Lås upp Dokumentmall Räntefaktura")
    ,@(450, "Button", 24094, "Kunder", "Lås upp Dokumentmall Kontantnota")
    ,@(451, "Button", 23896, "Kunder", "Lås upp Dokumentmall Extra fakturadokument")
    ,@(420, "Button", 24089, "Kunder", "Lås upp Utskriftsval Offert")
    ,@(421, "Button", 23892, "Kunder", "Lås upp Utskriftsval Order")
    ,@(422, "Button", 24095, "Kunder", "Lås upp Utskriftsval Följesedel")
    ,@(423, "Button", 24174, "Kunder", "Lås upp Utskriftsval Extra orderdokument")
    ,@(424, "Button", 23893, "Kunder", "Lås upp Utskriftsval Faktura")
    ,@(425, "Button", 23894, "Kunder", "Lås upp Utskriftsval Kreditfaktura")
    ,@(426, "Button", 23895, "Kunder", "Lås upp Utskriftsval Räntefaktura")
    ,@(427, "Button", 23984, "Kunder", "Lås upp Utskriftsval Kontantnota")
    ,@(428, "Button", 23985, "Kunder", "Lås upp Utskriftsval Extra fakturadokument")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$ws.Range("C451").Select()
